$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B11: hours worked this week, 15 -> 23
$ws.Range("B11").Value = 23

# Update C11: append note about correcting failures in Session 4 data
$ws.Range("C11").Value = "preprocessing: artifact rejection and ICA in session 3 and building pipeline session 4, working on exporting data" + [char]10 + "one session, correcting the failures in data of Session 4"

# Update the active selection to C19 (last thing the author clicked before saving)
$ws.Range("C19").Select()
